# Generate Report for Handoff
# - Status moves from "In Translation" to "Ready for handoff"
# - Handoff timestamps are refreshed to reflect the new handoff generation time
# - The "Status"/datetime columns grow wider to fit the new, longer text

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-26 18:41:22"

# --- zh-cn sheet ---
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-26 18:41:17"

# --- de-de sheet ---
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-26 18:41:22"

# --- Column width adjustments (status / datetime columns widened) ---
$wsOverview.Columns.Item(5).ColumnWidth = 16.3   # E: zh-cn status/date column
$wsOverview.Columns.Item(6).ColumnWidth = 16.3   # F: de-de status/date column
$wsZhCn.Columns.Item(3).ColumnWidth = 16.3        # C: Status column
$wsDeDe.Columns.Item(3).ColumnWidth = 16.3        # C: Status column
